$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add one test case: update the answer to "Any extra tests you would like to
# mention?" (row 19, column B) with a note about running the automatic test.
$ws.Range("B19").Value = "If you want to check our automatic test, do not run docker as a daemon since we're printing test information to stdout"

# The longer text needs a taller row to display, so bump row 19's height.
$ws.Rows.Item(19).RowHeight = 63

# Scroll the sheet so row 9 (column A) is the top-left visible cell, without
# changing the current selection (still B20).
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
